$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be auto-parsed as numbers
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "54.150.49"
$ws.Range("E2").Value = "  +5.11%  "
$ws.Range("D3").Value = "3.160.65"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "398.81"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").Value = "109.46"
$ws.Range("E6").Value = "  +6.32%  "
$ws.Range("D7").Value = "0.548"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.616"
$ws.Range("E9").Value = "  +5.39%  "
$ws.Range("D10").Value = "38.85"
$ws.Range("E10").Value = "  +5.69%  "
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "0.0878"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "3.640.43"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").Value = "19.08"
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("D15").Value = "8.03"
$ws.Range("E15").Value = "  +3.27%  "
$ws.Range("D16").Value = "1.06"
$ws.Range("E16").Value = "  +8.41%  "
$ws.Range("D17").Value = "3.161.42"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "10.49"
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").Value = "53.839.44"
$ws.Range("E19").Value = "  +4.30%  "
$ws.Range("D20").Value = "3.28"
$ws.Range("E20").Value = "  +4.43%  "
$ws.Range("D21").Value = "12.81"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("D23").Value = "71.04"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "271.61"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").Value = "3.25"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").Value = "7.99"
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").Value = "27.66"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("D28").Value = "7.40"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "0.169"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "0.111"
$ws.Range("E31").Value = "  +3.61%  "
$ws.Range("D32").Value = "11.00"
$ws.Range("E32").Value = "  +7.33%  "
$ws.Range("D33").Value = "0.0504"
$ws.Range("E33").Value = "  +12.56%  "
$ws.Range("D34").Value = "36.98"
$ws.Range("E34").Value = "  +6.50%  "
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("D36").Value = "50.49"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "3.65"
$ws.Range("E37").Value = "  +9.88%  "
$ws.Range("D38").Value = "0.996"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "2.82"
$ws.Range("E39").Value = "  +10.31%  "
$ws.Range("D40").Value = "4.10"
$ws.Range("E40").Value = "  +9.40%  "
$ws.Range("D41").Value = "0.290"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").Value = "17.31"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").Value = "1.90"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("D44").Value = "130.55"
$ws.Range("E44").Value = "  +4.23%  "
$ws.Range("D45").Value = "0.118"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("D46").Value = "22.20"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").Value = "2.080.53"
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").Value = "0.0342"
$ws.Range("E50").Value = "  +7.25%  "
$ws.Range("D51").Value = "0.0500"
$ws.Range("E51").Value = "  +15.13%  "
